# Update odds values on Sheet1 as per the FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10
$ws.Range("J10").Value = 2.38

# Row 12
$ws.Range("G12").Value  = 2.7
$ws.Range("I12").Value  = 3.05
$ws.Range("K12").Value  = 1.82
$ws.Range("L12").Value  = 3.7
$ws.Range("O12").Value  = 1.5
$ws.Range("P12").Value  = 2.27
$ws.Range("Q12").Value  = 2.42
$ws.Range("S12").Value  = 1.55
$ws.Range("T12").Value  = 2.15
$ws.Range("U12").Value  = 1.91
$ws.Range("V12").Value  = 1.7
$ws.Range("W12").Value  = 6.6
$ws.Range("X12").Value  = 12.5
$ws.Range("Y12").Value  = 10
$ws.Range("AC12").Value = 5.7
$ws.Range("AD12").Value = 5.1
$ws.Range("AJ12").Value = 11
$ws.Range("AK12").Value = 45
$ws.Range("AN12").Value = 4.4
$ws.Range("AT12").Value = 2.15
$ws.Range("AU12").Value = 6.9
$ws.Range("AW12").Value = 4.75
$ws.Range("AX12").Value = 18
$ws.Range("AY12").Value = 27
$ws.Range("AZ12").Value = 100
$ws.Range("BA12").Value = 150
$ws.Range("BB12").Value = 400
